$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.107.72'
$ws.Range("E2").Value = '  +4.54%  '

$ws.Range("D3").Value = '3.346.73'
$ws.Range("E3").Value = '  +9.41%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '257.42'
$ws.Range("E5").Value = '  +10.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '619.99'
$ws.Range("E6").Value = '  +2.56%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.17'
$ws.Range("E7").Value = '  +7.94%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.384'
$ws.Range("E8").Value = '  +3.51%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.998'
$ws.Range("E9").Value = '  -0.14%  '

$ws.Range("D10").Value = '3.342.64'
$ws.Range("E10").Value = '  +9.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.796'
$ws.Range("E11").Value = '  +0.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.199'
$ws.Range("E12").Value = '  +2.68%  '

$ws.Range("D13").Value = '97.893.16'
$ws.Range("E13").Value = '  +5.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.75'
$ws.Range("E14").Value = '  +8.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000245'
$ws.Range("E15").Value = '  +3.98%  '

$ws.Range("D16").Value = '3.962.62'
$ws.Range("E16").Value = '  +9.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.51'
$ws.Range("E17").Value = '  +5.83%  '

$ws.Range("D18").Value = '3.375.01'
$ws.Range("E18").Value = '  +10.15%  '

$ws.Range("E19").Value = '  +4.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.91'
$ws.Range("E20").Value = '  +4.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '484.13'
$ws.Range("E21").Value = '  +11.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.80'
$ws.Range("E22").Value = '  +3.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000204'
$ws.Range("E23").Value = '  +9.70%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.24'
$ws.Range("E24").Value = '  +6.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.66'
$ws.Range("E25").Value = '  +4.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.15'
$ws.Range("E26").Value = '  +4.74%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.07'
$ws.Range("E27").Value = '  +4.41%  '

$ws.Range("D28").Value = '3.556.20'
$ws.Range("E28").Value = '  +10.10%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("E30").Value = '  +5.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.244'
$ws.Range("E31").Value = '  -0.54%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.122'
$ws.Range("E32").Value = '  +3.20%  '

$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.24%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.26'
$ws.Range("E34").Value = '  +3.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.26'
$ws.Range("E35").Value = '  +8.56%  '

$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '517.84'
$ws.Range("E36").Value = '  +13.39%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.38'
$ws.Range("E37").Value = '  -2.60%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.150'
$ws.Range("E38").Value = '  -1.67%  '

$ws.Range("E39").Value = '  +4.51%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.87'
$ws.Range("E40").Value = '  +3.88%  '

$ws.Range("E41").Value = '  +3.23%  '

$ws.Range("E42").Value = '  +3.47%  '

$ws.Range("B43").Value = 'MantraDAO'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.55'
$ws.Range("E43").Value = '  -4.35%  '

$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.27'
$ws.Range("E44").Value = '  +6.87%  '

$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.780'
$ws.Range("E46").Value = '  +18.42%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '161.58'
$ws.Range("E47").Value = '  +1.06%  '

$ws.Range("E48").Value = '  +6.01%  '

$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.54'
$ws.Range("E49").Value = '  +4.17%  '

$ws.Range("B50").Value = 'ImmutableX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.37'
$ws.Range("E50").Value = '  +8.79%  '

$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.802'
$ws.Range("E51").Value = '  +13.39%  '
